$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl10"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.4153035
$ws.Range("H2").Value = 16.830607
$ws.Range("I2").Value = 0.06421843997858546
$ws.Range("J2").Value = 0.04492292295833487
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1626153333333333
$ws.Range("N2").Value = 0.487846
$ws.Range("O2").Value = 0.1293260700537641
$ws.Range("P2").Value = 0.1293260700537641
$ws.Range("Q2").Value = 1.368457383753667
$ws.Range("R2").Value = 8.210744302522
$ws.Range("S2").Value = 0.008305118467413991
$ws.Range("T2").Value = 0.005809705081529465

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl10"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.4153035
$ws.Range("H3").Value = 16.830607
$ws.Range("I3").Value = 0.06421843997858546
$ws.Range("J3").Value = 0.04492292295833487
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8767803333333334
$ws.Range("N3").Value = 2.630341
$ws.Range("O3").Value = 0.6972931302732585
$ws.Range("P3").Value = 0.6972931302732585
$ws.Range("Q3").Value = 7.378372607831167
$ws.Range("R3").Value = 44.270235646987
$ws.Range("S3").Value = 0.04477907703393322
$ws.Range("T3").Value = 0.03132444557064175

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl10"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.4153035
$ws.Range("H4").Value = 16.830607
$ws.Range("I4").Value = 0.06421843997858546
$ws.Range("J4").Value = 0.04492292295833487
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.21801
$ws.Range("N4").Value = 0.65403
$ws.Range("O4").Value = 0.1733807996729775
$ws.Range("P4").Value = 0.1733807996729775
$ws.Range("Q4").Value = 1.834620316035
$ws.Range("R4").Value = 11.00772189621
$ws.Range("S4").Value = 0.01113424447723825
$ws.Range("T4").Value = 0.007788772306163657

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl10"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.268178
$ws.Range("H5").Value = 75.80453399999999
$ws.Range("I5").Value = 0.192825246559582
$ws.Range("J5").Value = 0.202331457253709
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1626153333333333
$ws.Range("N5").Value = 0.487846
$ws.Range("O5").Value = 0.1293260700537641
$ws.Range("P5").Value = 0.1293260700537641
$ws.Range("Q5").Value = 4.108993188195999
$ws.Range("R5").Value = 36.980938693764
$ws.Range("S5").Value = 0.02493733134469885
$ws.Range("T5").Value = 0.02616673221487335

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl10"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 25.268178
$ws.Range("H6").Value = 75.80453399999999
$ws.Range("I6").Value = 0.192825246559582
$ws.Range("J6").Value = 0.202331457253709
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8767803333333334
$ws.Range("N6").Value = 2.630341
$ws.Range("O6").Value = 0.6972931302732585
$ws.Range("P6").Value = 0.6972931302732585
$ws.Range("Q6").Value = 22.154641529566
$ws.Range("R6").Value = 199.391773766094
$ws.Range("S6").Value = 0.1344557197692438
$ws.Range("T6").Value = 0.1410843351811887

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl10"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 25.268178
$ws.Range("H7").Value = 75.80453399999999
$ws.Range("I7").Value = 0.192825246559582
$ws.Range("J7").Value = 0.202331457253709
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.21801
$ws.Range("N7").Value = 0.65403
$ws.Range("O7").Value = 0.1733807996729775
$ws.Range("P7").Value = 0.1733807996729775
$ws.Range("Q7").Value = 5.508715485779999
$ws.Range("R7").Value = 49.57843937201999
$ws.Range("S7").Value = 0.03343219544563938
$ws.Range("T7").Value = 0.03508038985764692

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Cxcl10"
$ws.Range("C8").Value = "Ccr3"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 40.48325633333334
$ws.Range("H8").Value = 121.449769
$ws.Range("I8").Value = 0.3089337855705213
$ws.Range("J8").Value = 0.3241641027025684
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1626153333333333
$ws.Range("N8").Value = 0.487846
$ws.Range("O8").Value = 0.1293260700537641
$ws.Range("P8").Value = 0.1293260700537641
$ws.Range("Q8").Value = 6.583198223063778
$ws.Range("R8").Value = 59.248784007574
$ws.Range("S8").Value = 0.03995319239466778
$ws.Range("T8").Value = 0.04192286945502795

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Cxcl10"
$ws.Range("C9").Value = "Ccr3"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 40.48325633333334
$ws.Range("H9").Value = 121.449769
$ws.Range("I9").Value = 0.3089337855705213
$ws.Range("J9").Value = 0.3241641027025684
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8767803333333334
$ws.Range("N9").Value = 2.630341
$ws.Range("O9").Value = 0.6972931302732585
$ws.Range("P9").Value = 0.6972931302732585
$ws.Range("Q9").Value = 35.49492298235878
$ws.Range("R9").Value = 319.454306841229
$ws.Range("S9").Value = 0.2154174063876364
$ws.Range("T9").Value = 0.226037401895696

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Cxcl10"
$ws.Range("C10").Value = "Ccr3"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 40.48325633333334
$ws.Range("H10").Value = 121.449769
$ws.Range("I10").Value = 0.3089337855705213
$ws.Range("J10").Value = 0.3241641027025684
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.21801
$ws.Range("N10").Value = 0.65403
$ws.Range("O10").Value = 0.1733807996729775
$ws.Range("P10").Value = 0.1733807996729775
$ws.Range("Q10").Value = 8.825754713230001
$ws.Range("R10").Value = 79.43179241907
$ws.Range("S10").Value = 0.05356318678821712
$ws.Range("T10").Value = 0.0562038313518445

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Cxcl10"
$ws.Range("C11").Value = "Ccr3"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.0550545
$ws.Range("H11").Value = 20.110109
$ws.Range("I11").Value = 0.07673162517426207
$ws.Range("J11").Value = 0.05367630990912667
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1626153333333333
$ws.Range("N11").Value = 0.487846
$ws.Range("O11").Value = 0.1293260700537641
$ws.Range("P11").Value = 0.1293260700537641
$ws.Range("Q11").Value = 1.635106039202334
$ws.Range("R11").Value = 9.810636235214
$ws.Range("S11").Value = 0.00992339953262579
$ws.Range("T11").Value = 0.006941746215535271

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Cxcl10"
$ws.Range("C12").Value = "Ccr3"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.0550545
$ws.Range("H12").Value = 20.110109
$ws.Range("I12").Value = 0.07673162517426207
$ws.Range("J12").Value = 0.05367630990912667
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.8767803333333334
$ws.Range("N12").Value = 2.630341
$ws.Range("O12").Value = 0.6972931302732585
$ws.Range("P12").Value = 0.6972931302732585
$ws.Range("Q12").Value = 8.816074036194834
$ws.Range("R12").Value = 52.89644421716901
$ws.Range("S12").Value = 0.05350443510871556
$ws.Range("T12").Value = 0.03742812215805245

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Cxcl10"
$ws.Range("C13").Value = "Ccr3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.0550545
$ws.Range("H13").Value = 20.110109
$ws.Range("I13").Value = 0.07673162517426207
$ws.Range("J13").Value = 0.05367630990912667
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.21801
$ws.Range("N13").Value = 0.65403
$ws.Range("O13").Value = 0.1733807996729775
$ws.Range("P13").Value = 0.1733807996729775
$ws.Range("Q13").Value = 2.192102431545
$ws.Range("R13").Value = 13.15261458927
$ws.Range("S13").Value = 0.01330379053292073
$ws.Range("T13").Value = 0.009306441535538945

# Row 14
$ws.Range("A14").Value = "Neutrophils"
$ws.Range("B14").Value = "Cxcl10"
$ws.Range("C14").Value = "Ccr3"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 32.434321
$ws.Range("H14").Value = 97.302963
$ws.Range("I14").Value = 0.247511155882218
$ws.Range("J14").Value = 0.2597133609302807
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1626153333333333
$ws.Range("N14").Value = 0.487846
$ws.Range("O14").Value = 0.1293260700537641
$ws.Range("P14").Value = 0.1293260700537641
$ws.Range("Q14").Value = 5.274317920855334
$ws.Range("R14").Value = 47.468861287698
$ws.Range("S14").Value = 0.03200964508471186
$ws.Range("T14").Value = 0.03358770830956801

# Row 15
$ws.Range("A15").Value = "Neutrophils"
$ws.Range("B15").Value = "Cxcl10"
$ws.Range("C15").Value = "Ccr3"
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 32.434321
$ws.Range("H15").Value = 97.302963
$ws.Range("I15").Value = 0.247511155882218
$ws.Range("J15").Value = 0.2597133609302807
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.8767803333333334
$ws.Range("N15").Value = 2.630341
$ws.Range("O15").Value = 0.6972931302732585
$ws.Range("P15").Value = 0.6972931302732585
$ws.Range("Q15").Value = 28.43777477782034
$ws.Range("R15").Value = 255.939973000383
$ws.Range("S15").Value = 0.1725878286626642
$ws.Range("T15").Value = 0.181096342416864

# Row 16
$ws.Range("A16").Value = "Neutrophils"
$ws.Range("B16").Value = "Cxcl10"
$ws.Range("C16").Value = "Ccr3"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 32.434321
$ws.Range("H16").Value = 97.302963
$ws.Range("I16").Value = 0.247511155882218
$ws.Range("J16").Value = 0.2597133609302807
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.21801
$ws.Range("N16").Value = 0.65403
$ws.Range("O16").Value = 0.1733807996729775
$ws.Range("P16").Value = 0.1733807996729775
$ws.Range("Q16").Value = 7.071006321210001
$ws.Range("R16").Value = 63.63905689089
$ws.Range("S16").Value = 0.04291368213484193
$ws.Range("T16").Value = 0.04502931020384868

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Cxcl10"
$ws.Range("C17").Value = "Ccr3"
$ws.Range("D17").Value = "Inflammatory-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 14.38574166666667
$ws.Range("H17").Value = 43.157225
$ws.Range("I17").Value = 0.1097797468348313
$ws.Range("J17").Value = 0.1151918462459805
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1626153333333333
$ws.Range("N17").Value = 0.487846
$ws.Range("O17").Value = 0.1293260700537641
$ws.Range("P17").Value = 0.1293260700537641
$ws.Range("Q17").Value = 2.339342176372222
$ws.Range("R17").Value = 21.05407958735
$ws.Range("S17").Value = 0.01419738322964588
$ws.Range("T17").Value = 0.01489730877723011

# Row 18
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Cxcl10"
$ws.Range("C18").Value = "Ccr3"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 14.38574166666667
$ws.Range("H18").Value = 43.157225
$ws.Range("I18").Value = 0.1097797468348313
$ws.Range("J18").Value = 0.1151918462459805
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.8767803333333334
$ws.Range("N18").Value = 2.630341
$ws.Range("O18").Value = 0.6972931302732585
$ws.Range("P18").Value = 0.6972931302732585
$ws.Range("Q18").Value = 12.61313537374722
$ws.Range("R18").Value = 113.518218363725
$ws.Range("S18").Value = 0.07654866331106533
$ws.Range("T18").Value = 0.08032248305081566

# Row 19
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Cxcl10"
$ws.Range("C19").Value = "Ccr3"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 14.38574166666667
$ws.Range("H19").Value = 43.157225
$ws.Range("I19").Value = 0.1097797468348313
$ws.Range("J19").Value = 0.1151918462459805
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.21801
$ws.Range("N19").Value = 0.65403
$ws.Range("O19").Value = 0.1733807996729775
$ws.Range("P19").Value = 0.1733807996729775
$ws.Range("Q19").Value = 3.13623554075
$ws.Range("R19").Value = 28.22611986675
$ws.Range("S19").Value = 0.01903370029412006
$ws.Range("T19").Value = 0.01997205441793477
